$d = $word.ActiveDocument

# Change 1: "Finished the the requirements specification document..." ->
# "This month, I have finished the requirements specification document..."
$d.Content.Find.Execute("Finished the the requirements specification document", $true, $false, $false, $false, $false,
                         $true, 1, $false, "This month, I have finished the requirements specification document", 2)

# Change 2: remove the stray <w:lastRenderedPageBreak/> rendering artifact in front of
# "Supervisor Meeting" by re-writing the paragraph's text in place (same text, but Word
# regenerates the run without the rendering-only field).
$d.Paragraphs(20).Range.Text = "Supervisor Meeting"

# Change 3: "proved with good advise" -> "provided with good advise"
$d.Content.Find.Execute("suggested improvements and proved with good advise", $true, $false, $false, $false, $false,
                         $true, 1, $false, "suggested improvements and provided with good advise", 2)
